# Fund sheet ("基金受益憑證") gains a real header row plus the
# property_category / category / date / legislator_name / legislator_id /
# source_file / index metadata columns (I:O), matching the layout already
# used on the other sheets (汽車, 存款, 股票).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基金受益憑證")

# --- Header row (row 1): replace the stray duplicate data row with real
#     column headers ---
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "owner"
$ws.Cells.Item(1, 4).Value = "dealer"
$ws.Cells.Item(1, 5).Value = "quantity"
$ws.Cells.Item(1, 6).Value = "face_value"
$ws.Cells.Item(1, 7).Value = "currency"
$ws.Cells.Item(1, 8).Value = "total"
$ws.Cells.Item(1, 9).Value = "property_category"
$ws.Cells.Item(1, 10).Value = "category"
$ws.Cells.Item(1, 11).Value = "date"
$ws.Cells.Item(1, 12).Value = "legislator_name"
$ws.Cells.Item(1, 13).Value = "legislator_id"
$ws.Cells.Item(1, 14).Value = "source_file"
$ws.Cells.Item(1, 15).Value = "index"

# Match the bold / centered / bordered look already used for B1:H1 on the
# newly added header cells I1:O1.
$newHeader = $ws.Range("I1:O1")
$newHeader.Font.Bold = $true
$newHeader.HorizontalAlignment = -4108
$newHeader.VerticalAlignment = -4160
$newHeader.Borders.LineStyle = 1

# --- Data rows (2-4): append the metadata columns I:O ---
$indices = @(66, 67, 68)
for ($r = 2; $r -le 4; $r++) {
    $idx = $indices[$r - 2]
    $ws.Cells.Item($r, 9).Value = "fund"
    $ws.Cells.Item($r, 10).Value = "normal"
    # Force text storage so "2012-02-29" isn't auto-converted to a date serial.
    $ws.Cells.Item($r, 11).NumberFormat = "@"
    $ws.Cells.Item($r, 11).Value = "2012-02-29"
    $ws.Cells.Item($r, 12).Value = "林岱樺"
    $ws.Cells.Item($r, 13).Value = 904
    $ws.Cells.Item($r, 14).Value = "tmp3bff1"
    $ws.Cells.Item($r, 15).Value = $idx
}

$wb.Save()
